$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.02310156321232606, 0.3926748527515752, 0.2150774634025442, 0.463764448187379, 0.4806736408029301, 14),
    @(-0.01768369464496384, 0.5371479443576889, 0.4289931411788038, 0.6549756798376591, 0.6814717874811226, 13),
    @(0.04484032762995968, 0.3900917434466585, 0.2911825458959873, 0.5396133299835979, 0.5616584768642613, 12),
    @(0.04618976013043215, 0.3482013799665239, 0.2021242622245639, 0.4495823197419622, 0.4690307507104655, 11),
    @(0.01802609018820156, 0.396353456599608, 0.271167765387455, 0.5207377126610431, 0.5485767701364259, 10),
    @(-0.08230540368573748, 0.4030200371981603, 0.2441143047403328, 0.4940792494532966, 0.5167278209294698, 9),
    @(-0.03533701630595494, 0.5782167566939017, 0.5048288183583536, 0.7105130669863529, 0.7586304303617115, 8),
    @(0.1068317712072389, 0.4608479654118836, 0.2863969179506793, 0.5351606468628642, 0.5664049249844474, 7),
    @(-0.005335926578848192, 0.3488635011072845, 0.1472458661147809, 0.3837262906223405, 0.4203104481247061, 6),
    @(0.03223045541561795, 0.3832847050435998, 0.1801749238854154, 0.4244701684281422, 0.47320202032155, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
